$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '291.71'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.87%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '30.91'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.49%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.875'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.27%'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.73%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.228'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '22.44%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.680'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.14%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.704'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.52%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8967'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.14%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1669'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.79%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07929'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4.90%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08115'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.48%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03095'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '3.58%'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.44%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001501'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.69%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005840'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2.38%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.478'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.55%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.078'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.11%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3324'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.60%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1298'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.57%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.026'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-7.94%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2096'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4.82%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04515'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.84%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001209'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.43%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004664'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '15.40%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001300'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '3.91%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003388'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01576'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-4.73%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04391'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.09%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007308'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.48%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009683'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.36%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002021'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-3.36%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009357'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-15.98%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00005729'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-4.38%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000749'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.23%'
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002893'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '20.49%'
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.241'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '4.95%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002096'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.23%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001997'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.23%'
